# A new form submission (row 51) was collected and appended to the first
# worksheet ("八位序列号收集收集结果yd5") of the results workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Submitter name
$ws.Range("A51").Value = "彼世节拍"

# Submission timestamp - same date/time serial + number format as the
# other rows in column B.
$ws.Range("B51").Value = 45921.188587963
$ws.Range("B51").NumberFormat = $ws.Range("B50").NumberFormat

# Serial number (required field)
$ws.Range("C51").Value = "fdef3cc8"

# QQ number - stored as text, like the rest of column D, even though it
# looks numeric.
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1762488480"
